$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove old / now-unused cells (frees shared-string slot for "amps") ---
$ws.Range("J9").ClearContents()
$ws.Range("J12").ClearContents()
$ws.Range("J13").ClearContents()
$ws.Range("J14").ClearContents()
$ws.Range("J15").ClearContents()
$ws.Range("L16").ClearContents()

# Clear old scratch area (rows 22-24) before rebuilding it
$ws.Range("A22:E24").ClearContents()
$ws.Range("H22:J24").ClearContents()

# --- New row 19: Div/(A*100) ---
$ws.Range("A19").Value = "Div/(A*100)"
$ws.Range("B19").Formula = "=B15*100"
$ws.Range("C19").Formula = "=C15*100"
$ws.Range("D19").Formula = "=D15*100"
$ws.Range("F19").Formula = "=F15*100"
$ws.Range("G19").Formula = "=G15*100"
$ws.Range("H19").Formula = "=H15*100"

# --- New row 20: Count ---
$ws.Range("A20").Value = "Count"
$ws.Range("B20").Value = 60
$ws.Range("C20").Value = 120
$ws.Range("D20").Value = 200

# --- New row 21: Count * 100000 ---
$ws.Range("A21").Value = "Count * 100000"
$ws.Range("B21").Formula = "=B20*100000"
$ws.Range("C21").Formula = "=C20*100000"
$ws.Range("D21").Formula = "=D20*100000"

# --- Row 22: mA ---
$ws.Range("A22").Value = "mA"
$ws.Range("B22").Formula = "=B21/B19"
$ws.Range("C22").Formula = "=C21/C19"
$ws.Range("D22").Formula = "=D21/D19"

# --- Sheet view / formatting tweaks ---
$ws.Range("A1").Select()
$ws.Range("B26").Select()
